$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.443.82"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +8.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.679.78"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.53"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3710"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3445"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.40"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +14.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.181"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07267"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.42"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.150"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.751"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.678.32"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001110"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9991"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06719"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "81.34"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.46"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.108"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.05"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.439.23"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +8.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.446"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.687"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.50"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.53"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.863.06"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.21"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.365"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.041"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9788"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08439"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.699"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.61"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06491"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.82%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.381"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.939"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02328"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.260"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2116"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6187"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9985"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.26"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5971"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.68%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.760"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.47"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.032"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07225"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +6.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "76.05"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.80%  "
